$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 corresponds to the 517992c0-b950-4758-a718-d41ffd9e7b11 handback entry.
# Column L = "Latest Handback DateTime", column R = "Error Detail".

$wsZhCn.Range("L4").Value = "2017-02-09 16:43:10"
$wsZhCn.Range("R4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3cc8b70e6d10bb3896f7808b950214e436bd217d/e2e/517992c0-b950-4758-a718-d41ffd9e7b11.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22696147c9e27aa276e40797d7d66f4a3ddd9c40/e2e/517992c0-b950-4758-a718-d41ffd9e7b11.md."

$wsDeDe.Range("L4").Value = "2017-02-09 16:43:34"
$wsDeDe.Range("R4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3cc8b70e6d10bb3896f7808b950214e436bd217d/e2e/517992c0-b950-4758-a718-d41ffd9e7b11.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22696147c9e27aa276e40797d7d66f4a3ddd9c40/e2e/517992c0-b950-4758-a718-d41ffd9e7b11.md."
